$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: update B2 total, move selection to C7
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Select()
$ws1.Range("B2").Value = 3205797.3900000006
$ws1.Range("C7").Select()

# ---------------------------------------------------------------------------
# Sheet2: refresh column B figures, move selection to E9
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Select()
$ws2.Range("B2").Value = 627679.12
$ws2.Range("B3").Value = 63170.08
$ws2.Range("B4").Value = 558319.12000000011
$ws2.Range("B5").Value = 394429.3
$ws2.Range("B6").Value = 286389.90999999992
$ws2.Range("B7").Value = 508812.29
$ws2.Range("B8").Value = 194455
$ws2.Range("E9").Select()

# ---------------------------------------------------------------------------
# Sheet3: refresh column B figures, move selection to E12, becomes active tab
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Select()
$ws3.Range("B2").Value = 329854.56
$ws3.Range("B3").Value = 492025.49999999971
$ws3.Range("B4").Value = 474854.6999999999
$ws3.Range("B5").Value = 470940.67
$ws3.Range("B6").Value = 605487.00999999978
$ws3.Range("B7").Value = 563723.9
$ws3.Range("E12").Select()

# ---------------------------------------------------------------------------
# Sheet5: shift the year series back by two years (2018-2025), append the two
# new trailing rows, move selection to D7
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Select()

$ws5.Range("A2").Value = 2018
$ws5.Range("B2").Value = 37416.74

$ws5.Range("A3").Value = 2019
$ws5.Range("B3").Value = 231494.31

$ws5.Range("A4").Value = 2020
$ws5.Range("B4").Value = 329854.56

$ws5.Range("A5").Value = 2021
$ws5.Range("B5").Value = 492025.49999999971

$ws5.Range("A6").Value = 2022
$ws5.Range("B6").Value = 474854.6999999999

$ws5.Range("A7").Value = 2023
$ws5.Range("B7").Value = 470940.67

$ws5.Range("A8").Value = 2024
$ws5.Range("B8").Value = 605487.00999999978

$ws5.Range("A9").Value = 2025
$ws5.Range("B9").Value = 563723.9

$ws5.Range("D7").Select()

# ---------------------------------------------------------------------------
# Final active tab is Sheet3 (activeTab index 2)
# ---------------------------------------------------------------------------
$ws3.Select()
$ws3.Range("E12").Select()
